# Recreate the "calculated rank and box office $ correlations" commit:
#  - add a new "Correlation" worksheet at the end of the workbook
#  - copy Rank / Box Office Gross values from "Weekly Data" (C101:D127)
#  - add a CORREL() formula comparing the two columns
#  - tidy up the "Weekly Data" sheet selection (it is no longer the active tab)

$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Data")

# --- add the new sheet after the last existing sheet -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Correlation"

# --- headers -------------------------------------------------------------
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "Box Office Gross"
$ws.Range("D1").Value = "Correlation:"

# --- data: Rank (A) and Box Office Gross (B), copied from Weekly Data C101:D127 --
$weekly.Range("C101:D127").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# box office gross column formatted as currency, matching "Weekly Data" column D
$ws.Range("B2:B28").NumberFormat = $weekly.Range("D101").NumberFormat
$ws.Columns.Item(2).ColumnWidth = 14

# --- correlation formula --------------------------------------------------
$ws.Range("D2").Formula = "=CORREL(A2:A28,B2:B28)"

# --- "Weekly Data" is no longer the active tab; update its selection -------
$weekly.Select()
$weekly.Range("C101:D127").Select()
$excel.ActiveWindow.ScrollRow = 93

# --- window geometry from the recorded workbookView -----------------------
$excel.ActiveWindow.Left = 4920
$excel.ActiveWindow.Top = 2420

# --- new "Correlation" sheet ends up the active / selected tab -------------
$ws.Select()
$ws.Range("D3").Select()
